$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(223).Insert()

$ws.Cells.Item(223, 1).Value = 10
$ws.Cells.Item(223, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(223, 3).Value = "La Araucanía"
$ws.Cells.Item(223, 4).Value = 44524
$ws.Cells.Item(223, 5).Value = 9
$ws.Cells.Item(223, 6).Value = 100114014
$ws.Cells.Item(223, 7).Value = "Betarraga"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 100
$ws.Cells.Item(223, 11).Value = 700
$ws.Cells.Item(223, 12).Value = 700
$ws.Cells.Item(223, 13).Value = 700
$ws.Cells.Item(223, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(223, 15).Value = "Región Metropolitana"
$ws.Cells.Item(223, 16).Value = 140
$ws.Cells.Item(223, 17).Value = 5
$ws.Cells.Item(223, 18).Value = "Hortaliza"
